$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before row 3, shifting existing rows 3-5 down to 5-7.
$ws.Rows("3:4").Insert()

# Row 3: CHEBI:35470 / central nervous system drug
$ws.Range("A3").Value = "CHEBI:35470"
$ws.Range("B3").Value = "central nervous system drug"
$ws.Range("C3").Value = "A class of drugs producing both physiological and psychological effects through a variety of mechanisms involving the central nervous system."
$ws.Range("D3").Value = "function"
$ws.Range("S3").Value = "Proposed"
$ws.Range("V3").Value = "BG"

# Row 4: CMO:0000000 / clinical measurement
$ws.Range("A4").Value = "CMO:0000000"
$ws.Range("B4").Value = "clinical measurement"
$ws.Range("C4").Value = "A quantitative or qualitative value which is the result of an act of assessing a morphological or physiological state or property in a single individual or sample or a group of individuals or samples, based on direct observation or experimental manipulation."
$ws.Range("D4").Value = "planned process"
$ws.Range("S4").Value = "Proposed"
$ws.Range("V4").Value = "BG"

$ws.Range("A1").Select()
